$wb = $excel.ActiveWorkbook

# --- DatosCuenta sheet ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokPreProdJuneOne"
$wsCuenta.Range("B2").Value = "SmokeNamePreProdOne"
$wsCuenta.Range("C2").Value = 27100128
$wsCuenta.Range("D2").Value = 128

# --- DatosHogar sheet ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 648

# --- DatosMotor sheet ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP030"
$wsMotor.Range("B2").Value = "ABC12SSMP030"
$wsMotor.Range("C2").Value = "ZAZ123SSMP030"
$wsMotor.Range("A2:C2").Select()

# --- DatosAP sheet ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Range("A2").Value = 21200129

# Restore original active sheet/tab (DatosAP was active before edits)
$wsAP.Activate()
